$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, otherwise Excel auto-converts them to numeric values and
# the exact original string formatting (e.g. "1.00", "0.300") is lost.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.855.67"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.812.10"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "597.78"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "167.53"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "3.809.83"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "6.31"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "36.21"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "4.448.25"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "3.805.88"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("E17").Value = "  +4.20%  "
$ws.Range("D18").Value = "67.845.14"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "7.13"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D21").Value = "461.93"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("D23").Value = "0.703"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "12.13"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").Value = "2.11"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "10.02"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "3.955.84"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "2.80"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  +5.09%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "29.80"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "9.11"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "3.44"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "0.996"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D44").Value = "48.16"
$ws.Range("E44").Value = "  +3.05%  "
$ws.Range("D45").Value = "43.77"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "0.300"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "150.41"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "8.35"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "397.63"
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").Value = "26.34"
$ws.Range("E51").Value = "  +4.47%  "
